$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split the "-> Écrire un programme ..." run so a leading
# space becomes its own run (same Helvetica/24pt formatting).
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("-> Écrire un programme", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $origStart = $rng.Start
    $ins = $d.Range($origStart, $origStart)
    $ins.InsertBefore(" ")
    # Force the newly inserted space to live in its own run by nudging a
    # character property on it (and back), which keeps the visible
    # formatting identical but prevents it from being merged back into
    # the following run. (Use the forward range [origStart, origStart+1)
    # rather than origStart-1, so this still works when origStart is 0.)
    $spaceRng = $d.Range($origStart, $origStart + 1)
    $spaceRng.Bold = 1
    $spaceRng.Bold = 0
}

# ---------------------------------------------------------------------
# Change 2: drop the stray "_GoBack" bookmark that currently sits right
# after the "VA" run (in "- TVA (à calculer si TVA existe)").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Changes 3 & 4: re-create "_GoBack" so it spans from the start of the
# "Écrivez votre note: " paragraph through the end of the word "note"
# (splitting the "note:" run into "note" + ":").
# ---------------------------------------------------------------------
$pr = $d.Content
$foundPara = $pr.Find.Execute("Écrivez votre note:", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($foundPara) {
    $paraStart = $pr.Start

    $nr = $d.Content
    $foundNote = $nr.Find.Execute("note:", $true, $false, $false, `
        $false, $false, $true, 1, $false, "", 0)
    if ($foundNote) {
        $noteEnd = $nr.Start + 4
        $bmRange = $d.Range($paraStart, $noteEnd)
        $newBookmark = $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}

Write-Host "edit.ps1 done"
